$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.503.86'
$ws.Range('E2').Value = '  -1.03%  '
$ws.Range('D3').Value = '2.989.60'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.36'
$ws.Range('E5').Value = '  +2.16%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.05'
$ws.Range('E6').Value = '  -0.83%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '2.985.57'
$ws.Range('E8').Value = '  -0.25%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.514'
$ws.Range('E9').Value = '  -0.92%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.148'
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.04'
$ws.Range('E11').Value = '  +4.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.454'
$ws.Range('E12').Value = '  +3.72%  '
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.26'
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.126'
$ws.Range('E15').Value = '  +2.55%  '
$ws.Range('D16').Value = '3.479.91'
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').Value = '61.448.57'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.89'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').Value = '2.993.97'
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '449.72'
$ws.Range('E20').Value = '  -1.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.02'
$ws.Range('E21').Value = '  +2.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.683'
$ws.Range('E22').Value = '  +1.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.32'
$ws.Range('E23').Value = '  -0.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.93'
$ws.Range('E24').Value = '  +1.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.60'
$ws.Range('E25').Value = '  +5.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.18'
$ws.Range('E26').Value = '  -2.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.99'
$ws.Range('E27').Value = '  -1.17%  '
$ws.Range('E28').Value = '  +0.08%  '
$ws.Range('E29').Value = '  +3.02%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.27'
$ws.Range('E31').Value = '  +2.39%  '
$ws.Range('E32').Value = '  -0.84%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.26'
$ws.Range('E33').Value = '  +1.70%  '
$ws.Range('E34').Value = '  +1.68%  '
$ws.Range('D35').Value = '0.0₃0827'
$ws.Range('E35').Value = '  +5.89%  '
$ws.Range('E36').Value = '  -0.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.77'
$ws.Range('E37').Value = '  +1.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.26'
$ws.Range('E38').Value = '  +0.57%  '
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.03'
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.87'
$ws.Range('E41').Value = '  -0.66%  '
$ws.Range('E42').Value = '  +6.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '387.81'
$ws.Range('E43').Value = '  +2.15%  '
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('E45').Value = '  -0.40%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '38.43'
$ws.Range('E46').Value = '  +1.13%  '
$ws.Range('D47').Value = '2.699.81'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.53'
$ws.Range('E48').Value = '  +2.78%  '
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.15'
$ws.Range('E51').Value = '  +1.06%  '
